# Update the "admin" credential row (row 4) on Sheet1:
#   A4: admin@gmail.com -> asdfgfsg   (hyperlink to mailto:admin@gmail.com is preserved)
#   B4: admin123         -> ssfsfs
# Rows 1-3 (headers / other credential rows) are left untouched.
# Finally, move the active selection to B4 (matching the saved selection state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "asdfgfsg"
$ws.Range("B4").Value = "ssfsfs"

$ws.Range("B4").Select()
